# Update extrapolation calibration results after removing < USD 5 price
# from the calibration (noise). The recalculated columns D..H (ABSM1_RN,
# M1_RN, CM2_RN, CMN3_RN, CMN4_RN) change for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3" = 115955.2255462724;  "E3" = -0.01670229194722833; "F3" = 0.1881691251702542; "G3" = -1.198330183539438; "H3" = 12.13031415497148
    "D4" = 116678.6140084049;  "E4" = -0.02515636061401031; "F4" = 0.2274818802112115; "G4" = -1.498601766421102; "H4" = 13.05208136015519
    "D5" = 117478.7415924875;  "E5" = -0.02689556455888626; "F5" = 0.2295997912101166; "G5" = -0.8488861550963592; "H5" = 8.5138113623125
    "D8" = 119719.4197123511;  "E8" = -0.05074716512198433; "F8" = 0.2240109318190357; "G8" = -0.9338520001640457; "H8" = 6.743017105960065
    "D9" = 121303.5659644731;  "E9" = -0.08198894648673169; "F9" = 0.3624737229804428; "G9" = -1.997861032331176; "H9" = 12.63015744621798
    "D10" = 122768.96845271;   "E10" = -0.1179712470468311; "F10" = 0.4482425263296893; "G10" = -1.924195159957787; "H10" = 9.697580848004085
    "D11" = 124776.9884208325; "E11" = -0.1959728421668885; "F11" = 0.7865192950407556; "G11" = -2.597082958437904; "H11" = 12.55898539888622
    "D12" = 115305.8864000354; "E12" = -0.002104883204642449; "F12" = 0.1325200838312176; "G12" = -0.654349094243695; "H12" = 6.020192798401761
    "D14" = 115289.3589912245; "E14" = -0.003987857671811959; "F14" = 0.129745548845275; "G14" = -0.6881634047849767; "H14" = 5.807982189863465
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
